$d = $word.ActiveDocument

# Locate the paragraph that ends with "...and that they are all still there once they have all crossed."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*and that they are all still there once they have all crossed.*") {
        $target = $p
    }
}

$r = $target.Range
$r.Collapse(0)            # collapse the range to its end (just before the paragraph mark)
$r.InsertParagraphAfter() # new paragraph inherits the same pPr/rPr formatting

# The newly created paragraph is the one immediately following $target
$newPara1 = $target.Next()
$newPara1.Range.InsertAfter("The constraints are that he cannot bring all 3 items at one time. That means that he cannot take just one trip to the other side, he has to go 3 times. Another constraint is that each time, he has to leave each item with another item that won" + [char]0x2019 + "t eat the other item. The sub- goal is that the cat doesn" + [char]0x2019 + "t eat the parrot, and the parrot doesn" + [char]0x2019 + "t eat the food.")

# Insert a second new (empty-ish) paragraph containing a single space, right after the first one
$r2 = $newPara1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.InsertAfter(" ")
